$d = $word.ActiveDocument

# Locate the paragraph containing the "LOQ4057..." requisito text, then
# remove the following three paragraphs:
#   1. the blank paragraph right after it
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "(c) 2020 . Contact: luizeleno@usp.br. ..." paragraph
# leaving the blank paragraph / page-break paragraph that follow intact.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4057: Opera*es Unit*rias III (Requisito fraco)*") {
        $target = $p
    }
}

$blank = $target.Next()
$jupiter = $blank.Next()
$copyright = $jupiter.Next()
$after = $copyright.Next()

$start = $blank.Range.Start
$end = $after.Range.Start

$r = $d.Range($start, $end)
$r.Delete()
